# Automatische test-sync: 2025-06-19 21:37:50
# Adds the new "Afmelding nieuwsbrief" log entry (row 20) to the Logs sheet,
# extends the conditional formatting ranges to include the new row, and
# refreshes the Dashboard category counts/ordering to reflect the new entry.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new incoming mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Afmelding nieuwsbrief"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("C20").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D20").Value = "Afmelding / Nieuwsbrief"
$logs.Range("F20").Value = "2025-06-19 21:37:13"
$logs.Range("G20").Value = "Nee"

# Extend the conditional formatting ranges (Categorie + Beantwoord columns)
# so they cover the newly added row as well.
$logs.Range("D2:D19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))
$logs.Range("G2:G19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))

# --- Dashboard sheet: recompute category totals / ordering ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B5").Value = 2
$dash.Range("A6").Value = "Offerte / Prijsaanvraag"
$dash.Range("B6").Value = 2
$dash.Range("A7").Value = "Openingstijden / Locatie"
$dash.Range("B7").Value = 1
$dash.Range("A8").Value = "Factuur / Administratie"
$dash.Range("B8").Value = 1
$dash.Range("A9").Value = "Sollicitatie / Vacature"
$dash.Range("B9").Value = 1
$dash.Range("A10").Value = "Klacht / Probleem"
$dash.Range("B10").Value = 1
